$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, $Value)
    # Force literal text so Excel does not reinterpret numeric-looking
    # strings (e.g. "1.00", "0.0490") as numbers and strip formatting.
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.Style = "Normal"
}

$ws.Range("D2").Value = "73.095.97"
$ws.Range("E2").Value = "  +2.87%  "
$ws.Range("D3").Value = "3.991.83"
$ws.Range("E3").Value = "  +1.11%  "
Set-TextValue $ws.Range("D4") "1.00"
$ws.Range("E4").Value = "  +0.05%  "
Set-TextValue $ws.Range("D5") "595.98"
$ws.Range("E5").Value = "  +11.13%  "
Set-TextValue $ws.Range("D6") "163.67"
$ws.Range("E6").Value = "  +10.67%  "
$ws.Range("E7").Value = "  -0.49%  "
Set-TextValue $ws.Range("D8") "0.999"
$ws.Range("E8").Value = "  -0.06%  "
Set-TextValue $ws.Range("D9") "0.750"
$ws.Range("E9").Value = "  +1.58%  "
$ws.Range("E10").Value = "  +2.22%  "
Set-TextValue $ws.Range("D11") "54.59"
$ws.Range("E11").Value = "  -0.99%  "
Set-TextValue $ws.Range("D12") "0.0000320"
$ws.Range("E12").Value = "  +0.99%  "
Set-TextValue $ws.Range("D13") "10.99"
$ws.Range("E13").Value = "  +3.66%  "
$ws.Range("D14").Value = "4.632.51"
$ws.Range("E14").Value = "  +1.37%  "
$ws.Range("D15").Value = "3.995.10"
$ws.Range("E15").Value = "  +1.28%  "
$ws.Range("E16").Value = "  +8.84%  "
Set-TextValue $ws.Range("D17") "14.15"
$ws.Range("E17").Value = "  +1.90%  "
Set-TextValue $ws.Range("D18") "20.44"
$ws.Range("E18").Value = "  -0.25%  "
$ws.Range("E19").Value = "  +0.54%  "
$ws.Range("D20").Value = "72.833.92"
$ws.Range("E20").Value = "  +2.72%  "
Set-TextValue $ws.Range("D21") "437.73"
$ws.Range("E21").Value = "  +3.91%  "
Set-TextValue $ws.Range("D22") "4.74"
$ws.Range("E22").Value = "  +12.28%  "
Set-TextValue $ws.Range("D23") "96.36"
$ws.Range("E23").Value = "  -1.01%  "
$ws.Range("E24").Value = "  -4.40%  "
Set-TextValue $ws.Range("D25") "14.36"
$ws.Range("E25").Value = "  -0.63%  "
Set-TextValue $ws.Range("D26") "4.34"
$ws.Range("E26").Value = "  +13.77%  "
Set-TextValue $ws.Range("D27") "11.34"
$ws.Range("E27").Value = "  -0.39%  "
Set-TextValue $ws.Range("D28") "5.96"
$ws.Range("E28").Value = "  +1.11%  "
Set-TextValue $ws.Range("D29") "10.38"
$ws.Range("E29").Value = "  -2.77%  "
Set-TextValue $ws.Range("D30") "36.34"
$ws.Range("E30").Value = "  -0.23%  "
Set-TextValue $ws.Range("D31") "7.81"
$ws.Range("E31").Value = "  -1.28%  "
Set-TextValue $ws.Range("D32") "13.84"
$ws.Range("E32").Value = "  +3.73%  "
$ws.Range("E33").Value = "  +0.06%  "
$ws.Range("B34").Value = "Bittensor"
$ws.Range("C34").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws.Range("D34") "673.67"
$ws.Range("E34").Value = "  -1.49%  "
$ws.Range("B35").Value = "InjectiveProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D35") "48.28"
$ws.Range("E35").Value = "  -4.93%  "
Set-TextValue $ws.Range("D36") "70.82"
$ws.Range("E36").Value = "  +8.23%  "
$ws.Range("E37").Value = "  +10.86%  "
$ws.Range("E38").Value = "  -1.04%  "
Set-TextValue $ws.Range("D39") "3.39"
$ws.Range("E39").Value = "  +0.53%  "
$ws.Range("E40").Value = "  +0.10%  "
$ws.Range("B41").Value = "WEMIXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D41") "3.35"
$ws.Range("E41").Value = "  +5.66%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws.Range("D42") "0.145"
$ws.Range("E42").Value = "  -1.83%  "
Set-TextValue $ws.Range("D43") "1.00"
$ws.Range("E43").Value = "  +0.35%  "
Set-TextValue $ws.Range("D44") "0.0490"
$ws.Range("E44").Value = "  +1.85%  "
Set-TextValue $ws.Range("D45") "10.64"
$ws.Range("E45").Value = "  +6.35%  "
$ws.Range("E46").Value = "  +0.51%  "
$ws.Range("B47").Value = "Fetch.AI"
$ws.Range("C47").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue $ws.Range("D47") "2.64"
$ws.Range("E47").Value = "  -1.02%  "
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextValue $ws.Range("D48") "3.42"
$ws.Range("E48").Value = "  +2.16%  "
$ws.Range("D49").Value = "2.900.59"
$ws.Range("E49").Value = "  +10.39%  "
$ws.Range("E50").Value = "  +2.01%  "
Set-TextValue $ws.Range("D51") "3.40"
$ws.Range("E51").Value = "  +4.46%  "
